# Applies the timetable reshuffle described in the commit:
# "favorised tp1 with tp2 than tp1 with td2 (and the opposite)"
#
# Strategy for the simple "DPR110 communication" block moves (sheets G1-L1..G4-L1):
#   - The same text block moves from one empty-template cell to another within the
#     same sheet. We copy format (xlPasteFormats) + the text value to the new
#     location, then reset the old location back to the plain "empty" style (5)
#     and clear its content.
#
# Strategy for RSS-L2: full content/format rewrite of the B3:G7 block.

$xlPasteFormats = -4122

function Move-Block {
    param($ws, [string[]]$FromAddrs, [string[]]$ToAddrs)

    # Capture the text (identical across all source cells) and a style donor
    # cell before anything is modified.
    $text = $ws.Range($FromAddrs[0]).Value2
    $styleDonor = $FromAddrs[0]
    $emptyDonor = $null

    # find an already "empty" (style 5) cell on the same sheet to use as the
    # blank-style donor -- B3 is always an empty/plain cell on these sheets.
    $emptyDonor = "B3"

    foreach ($addr in $ToAddrs) {
        $ws.Range($styleDonor).Copy()
        $ws.Range($addr).PasteSpecial($xlPasteFormats)
        $ws.Range($addr).Value2 = $text
    }

    foreach ($addr in $FromAddrs) {
        $ws.Range($emptyDonor).Copy()
        $ws.Range($addr).PasteSpecial($xlPasteFormats)
        $ws.Range($addr).ClearContents()
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# G1-L1
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("G1-L1")
Move-Block $ws @("G5","B7") @("C3","G3")

# ---------------------------------------------------------------------------
# G2-L1
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("G2-L1")
Move-Block $ws @("D3","D7") @("C5","E6")

# ---------------------------------------------------------------------------
# G3-L1
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("G3-L1")
Move-Block $ws @("G3","C7") @("D7","F7")

# ---------------------------------------------------------------------------
# G4-L1 (D4 keeps its content; only F4 -> D5 moves)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("G4-L1")
Move-Block $ws @("F4") @("D5")

# ---------------------------------------------------------------------------
# RSS-L2 : full reshuffle of the B3:G7 grid
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("RSS-L2")

function Set-Cell {
    param($ws, [string]$Addr, [int]$StyleLevel, [string]$Text)

    # Style donors: level 5 = plain/empty, 6 = yellow (CM), 7 = blue (TD), 8 = green (TP)
    # These particular addresses are chosen because their style level is the same
    # both before and after the whole edit, so they remain valid donors no matter
    # which order the Set-Cell calls run in.
    switch ($StyleLevel) {
        5 { $donor = "B5" }
        6 { $donor = "C3" }
        7 { $donor = "B3" }
        8 { $donor = "D4" }
    }

    $ws.Range($donor).Copy()
    $ws.Range($Addr).PasteSpecial($xlPasteFormats)

    if ([string]::IsNullOrEmpty($Text)) {
        $ws.Range($Addr).ClearContents()
    } else {
        $ws.Range($Addr).Value2 = $Text
    }
}

# Row 3
Set-Cell $ws "B3" 7 "[RSS310] Reseaux Mobile`n(TD) - TD1`nProf: Aloun`nSalle: 101 /// [DAS311] RO`n(TD) - TD2`nProf: abderrahmane`nSalle: 102"
Set-Cell $ws "C3" 6 "[PAV312] Projet Integrateur`n(CM)`nProf: Encadreur`nSalle: 201"
Set-Cell $ws "D3" 6 "[PAV311] SD & Comp.Algo`n(CM)`nProf: Meyara`nSalle: 101"
Set-Cell $ws "E3" 5 $null
Set-Cell $ws "F3" 8 "[RSS311] Administration reseaux`n(TP) - TD1`nProf: Aloun`nSalle: 102 /// [DAS311] RO`n(TP) - TD2`nProf: abderrahmane`nSalle: 103"
Set-Cell $ws "G3" 7 "[DAS311] RO`n(TD) - TD1`nProf: abderrahmane`nSalle: 102"

# Row 4
Set-Cell $ws "B4" 6 "[DPR310] Communication`n(CM)`nProf: Dieynaba`nSalle: 101"
Set-Cell $ws "C4" 6 "[PAV310] POO JAVA`n(CM)`nProf: Esseyssah`nSalle: 101"
Set-Cell $ws "D4" 8 "[DAS311] RO`n(TP) - TD1`nProf: abderrahmane`nSalle: 102 /// [RSS311] Administration reseaux`n(TP) - TD2`nProf: Aloun`nSalle: 103"
Set-Cell $ws "E4" 5 $null
Set-Cell $ws "F4" 6 "[DAS311] RO`n(CM)`nProf: Cheikh`nSalle: 101"
Set-Cell $ws "G4" 5 $null

# Row 5
Set-Cell $ws "B5" 5 $null
Set-Cell $ws "C5" 6 "[DPR313] Gestion d'enterprise`n(CM)`nProf: El Bennany`nSalle: 201"
Set-Cell $ws "D5" 7 "[DAS311] RO`n(TD) - TD1`nProf: abderrahmane`nSalle: 102 /// [RSS310] Reseaux Mobile`n(TD) - TD2`nProf: Aloun`nSalle: 103"
Set-Cell $ws "E5" 7 "[DAS311] RO`n(TD) - TD2`nProf: abderrahmane`nSalle: 101"
Set-Cell $ws "F5" 5 $null
Set-Cell $ws "G5" 6 "[DPR310] Communication`n(CM)`nProf: Dieynaba`nSalle: 101"

# Row 6
Set-Cell $ws "B6" 7 "[DAS311] RO`n(TD) - TD1`nProf: abderrahmane`nSalle: 101 /// [RSS310] Reseaux Mobile`n(TD) - TD2`nProf: Aloun`nSalle: 102"
Set-Cell $ws "C6" 5 $null
Set-Cell $ws "D6" 6 "[RSS321] BD & CSI`n(CM)`nProf: Med Lemine`nSalle: 101"
Set-Cell $ws "E6" 6 "[DAS310] Maching Learning`n(CM)`nProf: Louly`nSalle: 201"
Set-Cell $ws "F6" 6 "[DPR311] Anglais`n(CM)`nProf: Blake`nSalle: 101"
Set-Cell $ws "G6" 6 "[PAV312] Projet Integrateur`n(CM)`nProf: Encadreur`nSalle: 101"

# Row 7
Set-Cell $ws "B7" 6 "[RSS311] Administration reseaux`n(CM)`nProf: Aloun`nSalle: 101"
Set-Cell $ws "C7" 7 "[RSS310] Reseaux Mobile`n(TD) - TD1`nProf: Aloun`nSalle: 101 /// [DAS311] RO`n(TD) - TD2`nProf: abderrahmane`nSalle: 102"
Set-Cell $ws "D7" 6 "[RSS321] BD & CSI`n(CM)`nProf: Med Lemine`nSalle: 201"
Set-Cell $ws "E7" 6 "[DPR311] Anglais`n(CM)`nProf: Blake`nSalle: 101"
Set-Cell $ws "F7" 6 "[RSS310] Reseaux Mobile`n(CM)`nProf: Aloun`nSalle: 201"
Set-Cell $ws "G7" 5 $null
